# Apply rotation of B/D/E/F/G values within each block of consecutive rows.
# Each row in a block takes on the previous (pre-edit) values of the next row
# in the block, with the last row wrapping around to the first.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Block 0: rows [154, 155, 156]
$v_B_154 = $ws.Cells.Item(154, 2).Value2
$v_B_155 = $ws.Cells.Item(155, 2).Value2
$v_B_156 = $ws.Cells.Item(156, 2).Value2
$v_D_154 = $ws.Cells.Item(154, 4).Value2
$v_D_155 = $ws.Cells.Item(155, 4).Value2
$v_D_156 = $ws.Cells.Item(156, 4).Value2
$v_E_154 = $ws.Cells.Item(154, 5).Value2
$v_E_155 = $ws.Cells.Item(155, 5).Value2
$v_E_156 = $ws.Cells.Item(156, 5).Value2
$v_F_154 = $ws.Cells.Item(154, 6).Value2
$v_F_155 = $ws.Cells.Item(155, 6).Value2
$v_F_156 = $ws.Cells.Item(156, 6).Value2
$v_G_154 = $ws.Cells.Item(154, 7).Value2
$v_G_155 = $ws.Cells.Item(155, 7).Value2
$v_G_156 = $ws.Cells.Item(156, 7).Value2
$ws.Cells.Item(154, 2).Value = $v_B_155
$ws.Cells.Item(155, 2).Value = $v_B_156
$ws.Cells.Item(156, 2).Value = $v_B_154
$ws.Cells.Item(154, 4).Value = $v_D_155
$ws.Cells.Item(155, 4).Value = $v_D_156
$ws.Cells.Item(156, 4).Value = $v_D_154
$ws.Cells.Item(154, 5).Value = $v_E_155
$ws.Cells.Item(155, 5).Value = $v_E_156
$ws.Cells.Item(156, 5).Value = $v_E_154
$ws.Cells.Item(154, 6).Value = $v_F_155
$ws.Cells.Item(155, 6).Value = $v_F_156
$ws.Cells.Item(156, 6).Value = $v_F_154
$ws.Cells.Item(154, 7).Value = $v_G_155
$ws.Cells.Item(155, 7).Value = $v_G_156
$ws.Cells.Item(156, 7).Value = $v_G_154

# Block 1: rows [271, 272]
$v_B_271 = $ws.Cells.Item(271, 2).Value2
$v_B_272 = $ws.Cells.Item(272, 2).Value2
$v_D_271 = $ws.Cells.Item(271, 4).Value2
$v_D_272 = $ws.Cells.Item(272, 4).Value2
$v_E_271 = $ws.Cells.Item(271, 5).Value2
$v_E_272 = $ws.Cells.Item(272, 5).Value2
$v_F_271 = $ws.Cells.Item(271, 6).Value2
$v_F_272 = $ws.Cells.Item(272, 6).Value2
$v_G_271 = $ws.Cells.Item(271, 7).Value2
$v_G_272 = $ws.Cells.Item(272, 7).Value2
$ws.Cells.Item(271, 2).Value = $v_B_272
$ws.Cells.Item(272, 2).Value = $v_B_271
$ws.Cells.Item(271, 4).Value = $v_D_272
$ws.Cells.Item(272, 4).Value = $v_D_271
$ws.Cells.Item(271, 5).Value = $v_E_272
$ws.Cells.Item(272, 5).Value = $v_E_271
$ws.Cells.Item(271, 6).Value = $v_F_272
$ws.Cells.Item(272, 6).Value = $v_F_271
$ws.Cells.Item(271, 7).Value = $v_G_272
$ws.Cells.Item(272, 7).Value = $v_G_271

# Block 2: rows [305, 306]
$v_B_305 = $ws.Cells.Item(305, 2).Value2
$v_B_306 = $ws.Cells.Item(306, 2).Value2
$v_D_305 = $ws.Cells.Item(305, 4).Value2
$v_D_306 = $ws.Cells.Item(306, 4).Value2
$v_E_305 = $ws.Cells.Item(305, 5).Value2
$v_E_306 = $ws.Cells.Item(306, 5).Value2
$v_F_305 = $ws.Cells.Item(305, 6).Value2
$v_F_306 = $ws.Cells.Item(306, 6).Value2
$v_G_305 = $ws.Cells.Item(305, 7).Value2
$v_G_306 = $ws.Cells.Item(306, 7).Value2
$ws.Cells.Item(305, 2).Value = $v_B_306
$ws.Cells.Item(306, 2).Value = $v_B_305
$ws.Cells.Item(305, 4).Value = $v_D_306
$ws.Cells.Item(306, 4).Value = $v_D_305
$ws.Cells.Item(305, 5).Value = $v_E_306
$ws.Cells.Item(306, 5).Value = $v_E_305
$ws.Cells.Item(305, 6).Value = $v_F_306
$ws.Cells.Item(306, 6).Value = $v_F_305
$ws.Cells.Item(305, 7).Value = $v_G_306
$ws.Cells.Item(306, 7).Value = $v_G_305

# Block 3: rows [308, 309]
$v_B_308 = $ws.Cells.Item(308, 2).Value2
$v_B_309 = $ws.Cells.Item(309, 2).Value2
$v_D_308 = $ws.Cells.Item(308, 4).Value2
$v_D_309 = $ws.Cells.Item(309, 4).Value2
$v_E_308 = $ws.Cells.Item(308, 5).Value2
$v_E_309 = $ws.Cells.Item(309, 5).Value2
$v_F_308 = $ws.Cells.Item(308, 6).Value2
$v_F_309 = $ws.Cells.Item(309, 6).Value2
$v_G_308 = $ws.Cells.Item(308, 7).Value2
$v_G_309 = $ws.Cells.Item(309, 7).Value2
$ws.Cells.Item(308, 2).Value = $v_B_309
$ws.Cells.Item(309, 2).Value = $v_B_308
$ws.Cells.Item(308, 4).Value = $v_D_309
$ws.Cells.Item(309, 4).Value = $v_D_308
$ws.Cells.Item(308, 5).Value = $v_E_309
$ws.Cells.Item(309, 5).Value = $v_E_308
$ws.Cells.Item(308, 6).Value = $v_F_309
$ws.Cells.Item(309, 6).Value = $v_F_308
$ws.Cells.Item(308, 7).Value = $v_G_309
$ws.Cells.Item(309, 7).Value = $v_G_308

# Block 4: rows [343, 344]
$v_B_343 = $ws.Cells.Item(343, 2).Value2
$v_B_344 = $ws.Cells.Item(344, 2).Value2
$v_D_343 = $ws.Cells.Item(343, 4).Value2
$v_D_344 = $ws.Cells.Item(344, 4).Value2
$v_E_343 = $ws.Cells.Item(343, 5).Value2
$v_E_344 = $ws.Cells.Item(344, 5).Value2
$v_F_343 = $ws.Cells.Item(343, 6).Value2
$v_F_344 = $ws.Cells.Item(344, 6).Value2
$v_G_343 = $ws.Cells.Item(343, 7).Value2
$v_G_344 = $ws.Cells.Item(344, 7).Value2
$ws.Cells.Item(343, 2).Value = $v_B_344
$ws.Cells.Item(344, 2).Value = $v_B_343
$ws.Cells.Item(343, 4).Value = $v_D_344
$ws.Cells.Item(344, 4).Value = $v_D_343
$ws.Cells.Item(343, 5).Value = $v_E_344
$ws.Cells.Item(344, 5).Value = $v_E_343
$ws.Cells.Item(343, 6).Value = $v_F_344
$ws.Cells.Item(344, 6).Value = $v_F_343
$ws.Cells.Item(343, 7).Value = $v_G_344
$ws.Cells.Item(344, 7).Value = $v_G_343

# Block 5: rows [347, 348]
$v_B_347 = $ws.Cells.Item(347, 2).Value2
$v_B_348 = $ws.Cells.Item(348, 2).Value2
$v_D_347 = $ws.Cells.Item(347, 4).Value2
$v_D_348 = $ws.Cells.Item(348, 4).Value2
$v_E_347 = $ws.Cells.Item(347, 5).Value2
$v_E_348 = $ws.Cells.Item(348, 5).Value2
$v_F_347 = $ws.Cells.Item(347, 6).Value2
$v_F_348 = $ws.Cells.Item(348, 6).Value2
$v_G_347 = $ws.Cells.Item(347, 7).Value2
$v_G_348 = $ws.Cells.Item(348, 7).Value2
$ws.Cells.Item(347, 2).Value = $v_B_348
$ws.Cells.Item(348, 2).Value = $v_B_347
$ws.Cells.Item(347, 4).Value = $v_D_348
$ws.Cells.Item(348, 4).Value = $v_D_347
$ws.Cells.Item(347, 5).Value = $v_E_348
$ws.Cells.Item(348, 5).Value = $v_E_347
$ws.Cells.Item(347, 6).Value = $v_F_348
$ws.Cells.Item(348, 6).Value = $v_F_347
$ws.Cells.Item(347, 7).Value = $v_G_348
$ws.Cells.Item(348, 7).Value = $v_G_347

# Block 6: rows [381, 382]
$v_B_381 = $ws.Cells.Item(381, 2).Value2
$v_B_382 = $ws.Cells.Item(382, 2).Value2
$v_D_381 = $ws.Cells.Item(381, 4).Value2
$v_D_382 = $ws.Cells.Item(382, 4).Value2
$v_E_381 = $ws.Cells.Item(381, 5).Value2
$v_E_382 = $ws.Cells.Item(382, 5).Value2
$v_F_381 = $ws.Cells.Item(381, 6).Value2
$v_F_382 = $ws.Cells.Item(382, 6).Value2
$v_G_381 = $ws.Cells.Item(381, 7).Value2
$v_G_382 = $ws.Cells.Item(382, 7).Value2
$ws.Cells.Item(381, 2).Value = $v_B_382
$ws.Cells.Item(382, 2).Value = $v_B_381
$ws.Cells.Item(381, 4).Value = $v_D_382
$ws.Cells.Item(382, 4).Value = $v_D_381
$ws.Cells.Item(381, 5).Value = $v_E_382
$ws.Cells.Item(382, 5).Value = $v_E_381
$ws.Cells.Item(381, 6).Value = $v_F_382
$ws.Cells.Item(382, 6).Value = $v_F_381
$ws.Cells.Item(381, 7).Value = $v_G_382
$ws.Cells.Item(382, 7).Value = $v_G_381

# Block 7: rows [392, 393]
$v_B_392 = $ws.Cells.Item(392, 2).Value2
$v_B_393 = $ws.Cells.Item(393, 2).Value2
$v_D_392 = $ws.Cells.Item(392, 4).Value2
$v_D_393 = $ws.Cells.Item(393, 4).Value2
$v_E_392 = $ws.Cells.Item(392, 5).Value2
$v_E_393 = $ws.Cells.Item(393, 5).Value2
$v_F_392 = $ws.Cells.Item(392, 6).Value2
$v_F_393 = $ws.Cells.Item(393, 6).Value2
$v_G_392 = $ws.Cells.Item(392, 7).Value2
$v_G_393 = $ws.Cells.Item(393, 7).Value2
$ws.Cells.Item(392, 2).Value = $v_B_393
$ws.Cells.Item(393, 2).Value = $v_B_392
$ws.Cells.Item(392, 4).Value = $v_D_393
$ws.Cells.Item(393, 4).Value = $v_D_392
$ws.Cells.Item(392, 5).Value = $v_E_393
$ws.Cells.Item(393, 5).Value = $v_E_392
$ws.Cells.Item(392, 6).Value = $v_F_393
$ws.Cells.Item(393, 6).Value = $v_F_392
$ws.Cells.Item(392, 7).Value = $v_G_393
$ws.Cells.Item(393, 7).Value = $v_G_392

# Block 8: rows [413, 414]
$v_B_413 = $ws.Cells.Item(413, 2).Value2
$v_B_414 = $ws.Cells.Item(414, 2).Value2
$v_D_413 = $ws.Cells.Item(413, 4).Value2
$v_D_414 = $ws.Cells.Item(414, 4).Value2
$v_E_413 = $ws.Cells.Item(413, 5).Value2
$v_E_414 = $ws.Cells.Item(414, 5).Value2
$v_F_413 = $ws.Cells.Item(413, 6).Value2
$v_F_414 = $ws.Cells.Item(414, 6).Value2
$v_G_413 = $ws.Cells.Item(413, 7).Value2
$v_G_414 = $ws.Cells.Item(414, 7).Value2
$ws.Cells.Item(413, 2).Value = $v_B_414
$ws.Cells.Item(414, 2).Value = $v_B_413
$ws.Cells.Item(413, 4).Value = $v_D_414
$ws.Cells.Item(414, 4).Value = $v_D_413
$ws.Cells.Item(413, 5).Value = $v_E_414
$ws.Cells.Item(414, 5).Value = $v_E_413
$ws.Cells.Item(413, 6).Value = $v_F_414
$ws.Cells.Item(414, 6).Value = $v_F_413
$ws.Cells.Item(413, 7).Value = $v_G_414
$ws.Cells.Item(414, 7).Value = $v_G_413

# Block 9: rows [449, 450]
$v_B_449 = $ws.Cells.Item(449, 2).Value2
$v_B_450 = $ws.Cells.Item(450, 2).Value2
$v_D_449 = $ws.Cells.Item(449, 4).Value2
$v_D_450 = $ws.Cells.Item(450, 4).Value2
$v_E_449 = $ws.Cells.Item(449, 5).Value2
$v_E_450 = $ws.Cells.Item(450, 5).Value2
$v_F_449 = $ws.Cells.Item(449, 6).Value2
$v_F_450 = $ws.Cells.Item(450, 6).Value2
$v_G_449 = $ws.Cells.Item(449, 7).Value2
$v_G_450 = $ws.Cells.Item(450, 7).Value2
$ws.Cells.Item(449, 2).Value = $v_B_450
$ws.Cells.Item(450, 2).Value = $v_B_449
$ws.Cells.Item(449, 4).Value = $v_D_450
$ws.Cells.Item(450, 4).Value = $v_D_449
$ws.Cells.Item(449, 5).Value = $v_E_450
$ws.Cells.Item(450, 5).Value = $v_E_449
$ws.Cells.Item(449, 6).Value = $v_F_450
$ws.Cells.Item(450, 6).Value = $v_F_449
$ws.Cells.Item(449, 7).Value = $v_G_450
$ws.Cells.Item(450, 7).Value = $v_G_449

# Block 10: rows [571, 572]
$v_B_571 = $ws.Cells.Item(571, 2).Value2
$v_B_572 = $ws.Cells.Item(572, 2).Value2
$v_D_571 = $ws.Cells.Item(571, 4).Value2
$v_D_572 = $ws.Cells.Item(572, 4).Value2
$v_E_571 = $ws.Cells.Item(571, 5).Value2
$v_E_572 = $ws.Cells.Item(572, 5).Value2
$v_F_571 = $ws.Cells.Item(571, 6).Value2
$v_F_572 = $ws.Cells.Item(572, 6).Value2
$v_G_571 = $ws.Cells.Item(571, 7).Value2
$v_G_572 = $ws.Cells.Item(572, 7).Value2
$ws.Cells.Item(571, 2).Value = $v_B_572
$ws.Cells.Item(572, 2).Value = $v_B_571
$ws.Cells.Item(571, 4).Value = $v_D_572
$ws.Cells.Item(572, 4).Value = $v_D_571
$ws.Cells.Item(571, 5).Value = $v_E_572
$ws.Cells.Item(572, 5).Value = $v_E_571
$ws.Cells.Item(571, 6).Value = $v_F_572
$ws.Cells.Item(572, 6).Value = $v_F_571
$ws.Cells.Item(571, 7).Value = $v_G_572
$ws.Cells.Item(572, 7).Value = $v_G_571

# Block 11: rows [582, 583]
$v_B_582 = $ws.Cells.Item(582, 2).Value2
$v_B_583 = $ws.Cells.Item(583, 2).Value2
$v_D_582 = $ws.Cells.Item(582, 4).Value2
$v_D_583 = $ws.Cells.Item(583, 4).Value2
$v_E_582 = $ws.Cells.Item(582, 5).Value2
$v_E_583 = $ws.Cells.Item(583, 5).Value2
$v_F_582 = $ws.Cells.Item(582, 6).Value2
$v_F_583 = $ws.Cells.Item(583, 6).Value2
$v_G_582 = $ws.Cells.Item(582, 7).Value2
$v_G_583 = $ws.Cells.Item(583, 7).Value2
$ws.Cells.Item(582, 2).Value = $v_B_583
$ws.Cells.Item(583, 2).Value = $v_B_582
$ws.Cells.Item(582, 4).Value = $v_D_583
$ws.Cells.Item(583, 4).Value = $v_D_582
$ws.Cells.Item(582, 5).Value = $v_E_583
$ws.Cells.Item(583, 5).Value = $v_E_582
$ws.Cells.Item(582, 6).Value = $v_F_583
$ws.Cells.Item(583, 6).Value = $v_F_582
$ws.Cells.Item(582, 7).Value = $v_G_583
$ws.Cells.Item(583, 7).Value = $v_G_582

# Block 12: rows [585, 586]
$v_B_585 = $ws.Cells.Item(585, 2).Value2
$v_B_586 = $ws.Cells.Item(586, 2).Value2
$v_D_585 = $ws.Cells.Item(585, 4).Value2
$v_D_586 = $ws.Cells.Item(586, 4).Value2
$v_E_585 = $ws.Cells.Item(585, 5).Value2
$v_E_586 = $ws.Cells.Item(586, 5).Value2
$v_F_585 = $ws.Cells.Item(585, 6).Value2
$v_F_586 = $ws.Cells.Item(586, 6).Value2
$v_G_585 = $ws.Cells.Item(585, 7).Value2
$v_G_586 = $ws.Cells.Item(586, 7).Value2
$ws.Cells.Item(585, 2).Value = $v_B_586
$ws.Cells.Item(586, 2).Value = $v_B_585
$ws.Cells.Item(585, 4).Value = $v_D_586
$ws.Cells.Item(586, 4).Value = $v_D_585
$ws.Cells.Item(585, 5).Value = $v_E_586
$ws.Cells.Item(586, 5).Value = $v_E_585
$ws.Cells.Item(585, 6).Value = $v_F_586
$ws.Cells.Item(586, 6).Value = $v_F_585
$ws.Cells.Item(585, 7).Value = $v_G_586
$ws.Cells.Item(586, 7).Value = $v_G_585

# Block 13: rows [591, 592]
$v_B_591 = $ws.Cells.Item(591, 2).Value2
$v_B_592 = $ws.Cells.Item(592, 2).Value2
$v_D_591 = $ws.Cells.Item(591, 4).Value2
$v_D_592 = $ws.Cells.Item(592, 4).Value2
$v_E_591 = $ws.Cells.Item(591, 5).Value2
$v_E_592 = $ws.Cells.Item(592, 5).Value2
$v_F_591 = $ws.Cells.Item(591, 6).Value2
$v_F_592 = $ws.Cells.Item(592, 6).Value2
$v_G_591 = $ws.Cells.Item(591, 7).Value2
$v_G_592 = $ws.Cells.Item(592, 7).Value2
$ws.Cells.Item(591, 2).Value = $v_B_592
$ws.Cells.Item(592, 2).Value = $v_B_591
$ws.Cells.Item(591, 4).Value = $v_D_592
$ws.Cells.Item(592, 4).Value = $v_D_591
$ws.Cells.Item(591, 5).Value = $v_E_592
$ws.Cells.Item(592, 5).Value = $v_E_591
$ws.Cells.Item(591, 6).Value = $v_F_592
$ws.Cells.Item(592, 6).Value = $v_F_591
$ws.Cells.Item(591, 7).Value = $v_G_592
$ws.Cells.Item(592, 7).Value = $v_G_591

# Block 14: rows [596, 597]
$v_B_596 = $ws.Cells.Item(596, 2).Value2
$v_B_597 = $ws.Cells.Item(597, 2).Value2
$v_D_596 = $ws.Cells.Item(596, 4).Value2
$v_D_597 = $ws.Cells.Item(597, 4).Value2
$v_E_596 = $ws.Cells.Item(596, 5).Value2
$v_E_597 = $ws.Cells.Item(597, 5).Value2
$v_F_596 = $ws.Cells.Item(596, 6).Value2
$v_F_597 = $ws.Cells.Item(597, 6).Value2
$v_G_596 = $ws.Cells.Item(596, 7).Value2
$v_G_597 = $ws.Cells.Item(597, 7).Value2
$ws.Cells.Item(596, 2).Value = $v_B_597
$ws.Cells.Item(597, 2).Value = $v_B_596
$ws.Cells.Item(596, 4).Value = $v_D_597
$ws.Cells.Item(597, 4).Value = $v_D_596
$ws.Cells.Item(596, 5).Value = $v_E_597
$ws.Cells.Item(597, 5).Value = $v_E_596
$ws.Cells.Item(596, 6).Value = $v_F_597
$ws.Cells.Item(597, 6).Value = $v_F_596
$ws.Cells.Item(596, 7).Value = $v_G_597
$ws.Cells.Item(597, 7).Value = $v_G_596

# Block 15: rows [701, 702]
$v_B_701 = $ws.Cells.Item(701, 2).Value2
$v_B_702 = $ws.Cells.Item(702, 2).Value2
$v_D_701 = $ws.Cells.Item(701, 4).Value2
$v_D_702 = $ws.Cells.Item(702, 4).Value2
$v_E_701 = $ws.Cells.Item(701, 5).Value2
$v_E_702 = $ws.Cells.Item(702, 5).Value2
$v_F_701 = $ws.Cells.Item(701, 6).Value2
$v_F_702 = $ws.Cells.Item(702, 6).Value2
$v_G_701 = $ws.Cells.Item(701, 7).Value2
$v_G_702 = $ws.Cells.Item(702, 7).Value2
$ws.Cells.Item(701, 2).Value = $v_B_702
$ws.Cells.Item(702, 2).Value = $v_B_701
$ws.Cells.Item(701, 4).Value = $v_D_702
$ws.Cells.Item(702, 4).Value = $v_D_701
$ws.Cells.Item(701, 5).Value = $v_E_702
$ws.Cells.Item(702, 5).Value = $v_E_701
$ws.Cells.Item(701, 6).Value = $v_F_702
$ws.Cells.Item(702, 6).Value = $v_F_701
$ws.Cells.Item(701, 7).Value = $v_G_702
$ws.Cells.Item(702, 7).Value = $v_G_701

# Block 16: rows [712, 713]
$v_B_712 = $ws.Cells.Item(712, 2).Value2
$v_B_713 = $ws.Cells.Item(713, 2).Value2
$v_D_712 = $ws.Cells.Item(712, 4).Value2
$v_D_713 = $ws.Cells.Item(713, 4).Value2
$v_E_712 = $ws.Cells.Item(712, 5).Value2
$v_E_713 = $ws.Cells.Item(713, 5).Value2
$v_F_712 = $ws.Cells.Item(712, 6).Value2
$v_F_713 = $ws.Cells.Item(713, 6).Value2
$v_G_712 = $ws.Cells.Item(712, 7).Value2
$v_G_713 = $ws.Cells.Item(713, 7).Value2
$ws.Cells.Item(712, 2).Value = $v_B_713
$ws.Cells.Item(713, 2).Value = $v_B_712
$ws.Cells.Item(712, 4).Value = $v_D_713
$ws.Cells.Item(713, 4).Value = $v_D_712
$ws.Cells.Item(712, 5).Value = $v_E_713
$ws.Cells.Item(713, 5).Value = $v_E_712
$ws.Cells.Item(712, 6).Value = $v_F_713
$ws.Cells.Item(713, 6).Value = $v_F_712
$ws.Cells.Item(712, 7).Value = $v_G_713
$ws.Cells.Item(713, 7).Value = $v_G_712

# Block 17: rows [864, 865]
$v_B_864 = $ws.Cells.Item(864, 2).Value2
$v_B_865 = $ws.Cells.Item(865, 2).Value2
$v_D_864 = $ws.Cells.Item(864, 4).Value2
$v_D_865 = $ws.Cells.Item(865, 4).Value2
$v_E_864 = $ws.Cells.Item(864, 5).Value2
$v_E_865 = $ws.Cells.Item(865, 5).Value2
$v_F_864 = $ws.Cells.Item(864, 6).Value2
$v_F_865 = $ws.Cells.Item(865, 6).Value2
$v_G_864 = $ws.Cells.Item(864, 7).Value2
$v_G_865 = $ws.Cells.Item(865, 7).Value2
$ws.Cells.Item(864, 2).Value = $v_B_865
$ws.Cells.Item(865, 2).Value = $v_B_864
$ws.Cells.Item(864, 4).Value = $v_D_865
$ws.Cells.Item(865, 4).Value = $v_D_864
$ws.Cells.Item(864, 5).Value = $v_E_865
$ws.Cells.Item(865, 5).Value = $v_E_864
$ws.Cells.Item(864, 6).Value = $v_F_865
$ws.Cells.Item(865, 6).Value = $v_F_864
$ws.Cells.Item(864, 7).Value = $v_G_865
$ws.Cells.Item(865, 7).Value = $v_G_864

